# Commit: "#5: property boat&car done"
# The 汽車 (car) sheet (sheet3 / Worksheets.Item(3)) is extended from a
# bare 6-column layout (name, capacity/cc, owner, register_date,
# register_reason, acquire_value) to the full common schema shared by the
# 土地 (land) and 建物 (building) sheets: a "capacity" header column is
# inserted, and property_category / category / date / legislator_name /
# legislator_id / source_file / index columns are appended, mirroring the
# row-2 data that already existed in the other property sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ---- Row 2 (data row) gets the new trailing columns --------------------
# (existing A2:G2 values are untouched)
$ws.Range("J2").NumberFormat = "@"   # keep "2012-03-06" literal text, not a date serial

$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2012-03-06"
$ws.Range("K2").Value = "吳育仁"
$ws.Range("L2").Value = 1734
$ws.Range("M2").Value = "tmp476d1"
$ws.Range("N2").Value = 38

# ---- Row 1 (header row) is rebuilt completely --------------------------
# Old header row was: B1=BMW C1=2500 D1=蔡瓊姿 E1=98年11月01曰 F1=貝買 G1=2250000
# (those were actually stray data values sitting in row 1). The new header
# row instead carries column names, with a new "capacity" column inserted
# right after "name".
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# ---- Formatting ----------------------------------------------------------
# Row 1 (B1:N1) uses the same bold/centered/bordered header style as the
# rest of the sheet. Copy the format from an already-styled header cell so
# the engine dedupes/reuses the existing xf instead of inventing a new one.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("D1:N1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
